$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "61.732.20"
$ws.Range("E2").Value = "  -1.68%  "

# Row 3
$ws.Range("D3").Value = "2.998.37"
$ws.Range("E3").Value = "  -1.02%  "

# Row 5
$ws.Range("D5").Value = "'597.81"
$ws.Range("E5").Value = "  +2.28%  "

# Row 6
$ws.Range("D6").Value = "'144.27"

# Row 7
$ws.Range("E7").Value = "  +0.04%  "

# Row 8
$ws.Range("E8").Value = "  -0.72%  "

# Row 9
$ws.Range("D9").Value = "2.998.72"
$ws.Range("E9").Value = "  -1.03%  "

# Row 10
$ws.Range("E10").Value = "  -2.60%  "

# Row 11
$ws.Range("D11").Value = "'5.90"
$ws.Range("E11").Value = "  +3.76%  "

# Row 12
$ws.Range("E12").Value = "  +4.11%  "

# Row 13
$ws.Range("D13").Value = "'0.0000229"
$ws.Range("E13").Value = "  -1.02%  "

# Row 14
$ws.Range("D14").Value = "'34.31"
$ws.Range("E14").Value = "  -3.18%  "

# Row 15
$ws.Range("E15").Value = "  +2.34%  "

# Row 16
$ws.Range("D16").Value = "3.492.58"
$ws.Range("E16").Value = "  -1.20%  "

# Row 17
$ws.Range("D17").Value = "'7.00"
$ws.Range("E17").Value = "  -1.12%  "

# Row 18
$ws.Range("D18").Value = "61.673.73"
$ws.Range("E18").Value = "  -1.78%  "

# Row 19
$ws.Range("D19").Value = "2.998.30"
$ws.Range("E19").Value = "  -1.00%  "

# Row 20
$ws.Range("D20").Value = "'455.35"
$ws.Range("E20").Value = "  -2.89%  "

# Row 21
$ws.Range("D21").Value = "'14.02"
$ws.Range("E21").Value = "  -0.29%  "

# Row 22
$ws.Range("D22").Value = "'0.687"
$ws.Range("E22").Value = "  -0.84%  "

# Row 23
$ws.Range("D23").Value = "'7.36"
$ws.Range("E23").Value = "  -0.93%  "

# Row 24
$ws.Range("D24").Value = "'82.34"
$ws.Range("E24").Value = "  +1.60%  "

# Row 25
$ws.Range("E25").Value = "  -7.17%  "

# Row 26
$ws.Range("D26").Value = "'12.16"
$ws.Range("E26").Value = "  -1.88%  "

# Row 27
$ws.Range("D27").Value = "'10.36"
$ws.Range("E27").Value = "  -1.06%  "

# Row 28
$ws.Range("E28").Value = "  +0.08%  "

# Row 29
$ws.Range("E29").Value = "  +1.92%  "

# Row 30
$ws.Range("E30").Value = "  +0.02%  "

# Row 31
$ws.Range("D31").Value = "'7.03"
$ws.Range("E31").Value = "  -3.33%  "

# Row 32
$ws.Range("E32").Value = "  -4.00%  "

# Row 33
$ws.Range("D33").Value = "'27.44"
$ws.Range("E33").Value = "  -0.10%  "

# Row 34
$ws.Range("E34").Value = "  +0.06%  "

# Row 35
$ws.Range("D35").Value = "0.0₃0820"
$ws.Range("E35").Value = "  +2.98%  "

# Row 36
$ws.Range("E36").Value = "  -1.94%  "

# Row 37
$ws.Range("D37").Value = "'5.75"
$ws.Range("E37").Value = "  -0.61%  "

# Row 39
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").Value = "'50.34"
$ws.Range("E39").Value = "  +0.10%  "

# Row 40
$ws.Range("B40").Value = "Cosmos"
$ws.Range("C40").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D40").Value = "'9.16"
$ws.Range("E40").Value = "  +1.67%  "

# Row 41
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "'2.89"
$ws.Range("E41").Value = "  -2.34%  "

# Row 42
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.122"
$ws.Range("E42").Value = "  +7.94%  "

# Row 43
$ws.Range("D43").Value = "'402.54"
$ws.Range("E43").Value = "  -5.05%  "

# Row 44
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0354"
$ws.Range("E44").Value = "  -0.68%  "

# Row 45
$ws.Range("B45").Value = "Arweave"
$ws.Range("C45").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D45").Value = "'39.00"
$ws.Range("E45").Value = "  +3.05%  "

# Row 46
$ws.Range("E46").Value = "  -5.31%  "

# Row 47
$ws.Range("D47").Value = "2.719.66"
$ws.Range("E47").Value = "  -3.14%  "

# Row 48
$ws.Range("D48").Value = "'133.73"
$ws.Range("E48").Value = "  +3.16%  "

# Row 49
$ws.Range("E49").Value = "  +0.16%  "

# Row 50
$ws.Range("D50").Value = "'0.108"
$ws.Range("E50").Value = "  -0.64%  "

# Row 51
$ws.Range("E51").Value = "  +1.34%  "
